$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.001.57'
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = '  +1.40%  '

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.052.35'
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = '  -2.05%  '

$ws.Range("E4").Value = '  +0.14%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.43'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -1.24%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.674'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +3.10%  '

$ws.Range("E7").Value = '  +0.00%  '

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.52'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  +15.24%  '

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.63'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  +1.73%  '

$ws.Range("E10").Value = '  +1.87%  '

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0782'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  +5.47%  '

$ws.Range("E12").Value = '  +6.22%  '

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.85'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +2.46%  '

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.353.25'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  -1.91%  '

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.813'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -1.16%  '

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  +3.04%  '

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.055.55'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -1.92%  '

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.936.31'
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  +1.21%  '

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0925'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  +12.07%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.62'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +0.12%  '

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.13'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  +7.39%  '

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.32'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +3.71%  '

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.52'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("E24").Value = '  -0.04%  '

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -1.09%  '

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.01'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  +0.04%  '

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.93'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  -1.56%  '

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.98'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -6.01%  '

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.98'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -0.12%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.124'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  +1.74%  '

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.56'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +2.91%  '

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0620'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("E33").Value = '  +5.71%  '

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.32'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  +6.49%  '

$ws.Range("E35").Value = '  +0.02%  '

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0870'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  -5.74%  '

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.26'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -3.77%  '

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  -5.21%  '

$ws.Range("E39").Value = '  +0.54%  '

$ws.Range("E40").Value = '  +21.25%  '

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.71'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  +11.35%  '

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0222'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  +0.59%  '

$ws.Range("E43").Value = '  -1.79%  '

$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '95.95'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  -1.10%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +2.00%  '

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.15'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +52.19%  '

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.38'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  +7.23%  '

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.289.92'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -2.74%  '

$ws.Range("E49").Value = '  +2.85%  '

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '13.00'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -53.81%  '

$ws.Range("E51").Value = '  +6.75%  '
